# "Generate Report for Handback"
# Fills in the Latest Target File / Latest Handback File / Latest Handback
# DateTime columns on the zh-cn and de-de status sheets, flips the Status
# column from "Ready for handoff" to "Handed back: in sync with en-US",
# and widens a few columns so the new long file names are readable.

$wb = $excel.ActiveWorkbook

$shZhCn = $wb.Worksheets.Item("zh-cn")
$shDeDe = $wb.Worksheets.Item("de-de")
$shOverview = $wb.Worksheets.Item("Overview")

$mdFileName   = "94c0e964-2323-4f37-ae93-3328c6e77f30.md"
$mdTargetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41575a505568da6c8733ba8bbc79ad2a9a161c22/e2e/94c0e964-2323-4f37-ae93-3328c6e77f30.md"
$zhHandback   = "94c0e964-2323-4f37-ae93-3328c6e77f30.273af830355ba05e8f841c180ec6c5e3ce13173e.zh-cn.xlf"
$deHandback   = "94c0e964-2323-4f37-ae93-3328c6e77f30.273af830355ba05e8f841c180ec6c5e3ce13173e.de-de.xlf"
$handbackStatus = "Handed back: in sync with en-US"

# Exact ColumnWidth input values that round-trip to the desired stored
# <col width="..."> values of this workbook (29.9777047293527 / 40).
$wideColWidth  = 29.9777047293527
$fullColWidth  = 235.0 / 6.0

# --- Status column (C) on both language sheets: "Ready for handoff" -> "Handed back: in sync with en-US"
$shZhCn.Range("C2").Value = $handbackStatus
$shDeDe.Range("C2").Value = $handbackStatus

# --- zh-cn sheet: Latest Target File (I2), Latest Handback File (J2), Latest Handback DateTime (K2)
$shZhCn.Range("I2").Value = $mdFileName
$shZhCn.Hyperlinks.Add($shZhCn.Range("I2"), $mdTargetUrl, $null, $null, $mdFileName)
$shZhCn.Range("J2").Value = $zhHandback
$shZhCn.Range("K2").Value = "2016-08-22 20:58:58"

# --- de-de sheet: Latest Target File (I2), Latest Handback File (J2), Latest Handback DateTime (K2)
$shDeDe.Range("I2").Value = $mdFileName
$shDeDe.Hyperlinks.Add($shDeDe.Range("I2"), $mdTargetUrl, $null, $null, $mdFileName)
$shDeDe.Range("J2").Value = $deHandback
$shDeDe.Range("K2").Value = "2016-08-22 20:59:13"

# --- Column widths: widen Overview!E:F, and zh-cn/de-de column C (Status) and I:J (new file-name columns)
$shOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$shOverview.Columns.Item(6).ColumnWidth = $wideColWidth

$shZhCn.Columns.Item(3).ColumnWidth = $wideColWidth
$shZhCn.Columns.Item(9).ColumnWidth = $fullColWidth
$shZhCn.Columns.Item(10).ColumnWidth = $fullColWidth

$shDeDe.Columns.Item(3).ColumnWidth = $wideColWidth
$shDeDe.Columns.Item(9).ColumnWidth = $fullColWidth
$shDeDe.Columns.Item(10).ColumnWidth = $fullColWidth
